# T Natarajan.xlsx - "complete!! -> scrapping whole ipl"
#
# The scraper added a leading "matchNo" column to every batter sheet and
# started stamping the sheet tab with the batter's own name instead of the
# generic "Sheet1". For this row: matchNo = "6th".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the tab from "Sheet1" to the batter's name.
$ws.Name = "T Natarajan"

# Insert a new column before A; this shifts every existing header/value one
# column to the right (teamName: A->B, batterName: B->C, ..., result: L->M)
# while keeping each value intact.
$ws.Range("A1").EntireColumn.Insert()

# Populate the new leading column.
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "6th"
